$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: OnlineAccountInformation  (A1:B2 -> A1:E2, add T.C (Azure)/T.C (Desc.)/Error)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("OnlineAccountInformation")
$ws1.Range("C1").Value = "T.C (Azure)"
$ws1.Range("D1").Value = "T.C (Desc.)"
$ws1.Range("E1").Value = "Error"
$ws1.Range("C2").Value = 119804

$rng1 = $ws1.Range("A1:E2")
$rng1.Borders.LineStyle = 1
$rng1.Borders.Weight = 2
$rng1.BorderAround(1, -4138)

$hdr1 = $ws1.Range("A1:E1")
$hdr1.Font.Bold = $true
$hdr1.Font.Size = 14
$ws1.Rows(1).RowHeight = 18.75
$ws1.Rows(2).RowHeight = 15.75

$ws1.Range("C2").HorizontalAlignment = -4131
$ws1.Range("A1:E1").Select()

# ---------------------------------------------------------------------------
# Sheet 2: AccountBalanceInformation (A1:A2 -> A1:D2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("AccountBalanceInformation")
$ws2.Range("B1").Value = "T.C (Azure)"
$ws2.Range("C1").Value = "T.C (Desc.)"
$ws2.Range("D1").Value = "Error"
$ws2.Range("B2").Value = 119804

$rngNew2 = $ws2.Range("B1:D2")
$rngNew2.Borders.LineStyle = 1
$rngNew2.Borders.Weight = 2
$full2 = $ws2.Range("A1:D2")
$full2.BorderAround(1, -4138)
$ws2.Range("A1").Borders.Item(9).LineStyle = 1
$ws2.Range("A1").Borders.Item(9).Weight = -4138
$ws2.Range("A2").Borders.Item(10).LineStyle = 1
$ws2.Range("A2").Borders.Item(10).Weight = -4138

$hdr2 = $ws2.Range("A1:D1")
$hdr2.Font.Bold = $true
$hdr2.Font.Size = 14
$ws2.Rows(1).RowHeight = 18.75
$ws2.Rows(2).RowHeight = 15.75

$ws2.Range("B2").HorizontalAlignment = -4131
$ws2.Range("A1:D1").Select()

# ---------------------------------------------------------------------------
# Sheet 3: STOReversed (A1:A2 -> A1:D2)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("STOReversed")
$ws3.Range("B1").Value = "T.C (Azure)"
$ws3.Range("C1").Value = "T.C (Desc.)"
$ws3.Range("D1").Value = "Error"
$ws3.Range("B2").Value = 119803

$rngNew3 = $ws3.Range("B1:D2")
$rngNew3.Borders.LineStyle = 1
$rngNew3.Borders.Weight = 2
$full3 = $ws3.Range("A1:D2")
$full3.BorderAround(1, -4138)
$ws3.Range("A1").Borders.Item(9).LineStyle = 1
$ws3.Range("A1").Borders.Item(9).Weight = -4138
$ws3.Range("A2").Borders.Item(10).LineStyle = 1
$ws3.Range("A2").Borders.Item(10).Weight = -4138

$hdr3 = $ws3.Range("A1:D1")
$hdr3.Font.Bold = $true
$hdr3.Font.Size = 14
$ws3.Rows(1).RowHeight = 18.75
$ws3.Rows(2).RowHeight = 15.75

$ws3.Range("B2").HorizontalAlignment = -4131
$ws3.Range("A1:D1").Select()

# ---------------------------------------------------------------------------
# Sheet 4: UnauthorizedList (A1:B2 -> A1:E2), only new cols C:E get boxed, blank
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("UnauthorizedList")

$rng4 = $ws4.Range("C1:E2")
$rng4.Borders.LineStyle = 1
$rng4.Borders.Weight = 2
$rng4.BorderAround(1, -4138)

$hdr4 = $ws4.Range("C1:E1")
$hdr4.Font.Bold = $true
$hdr4.Font.Size = 14
$ws4.Rows(1).RowHeight = 18.75

$ws4.Range("C2").HorizontalAlignment = -4131
$ws4.Range("C1").Select()
$ws4.Range("C1:E2").Select()

# ---------------------------------------------------------------------------
# Sheet 5: StandingOrderINP (A1:F2 -> A1:I2)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("StandingOrderINP")
$ws5.Range("G1").Value = "T.C (Azure)"
$ws5.Range("H1").Value = "T.C (Desc.)"
$ws5.Range("I1").Value = "Error"
$ws5.Range("G2").Value = 119800

$rng5 = $ws5.Range("A1:I2")
$rng5.Borders.LineStyle = 1
$rng5.Borders.Weight = 2
$rng5.BorderAround(1, -4138)

$hdr5 = $ws5.Range("A1:I1")
$hdr5.Font.Bold = $true
$hdr5.Font.Size = 14
$ws5.Rows(1).RowHeight = 18.75
$ws5.Rows(2).RowHeight = 15.75

$ws5.Range("G2").HorizontalAlignment = -4131
$ws5.Range("A1:I1").Select()

# ---------------------------------------------------------------------------
# Sheet 6: StandingOrderAuth (A1:F2 -> A1:I2)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("StandingOrderAuth")
$ws6.Range("G1").Value = "T.C (Azure)"
$ws6.Range("H1").Value = "T.C (Desc.)"
$ws6.Range("I1").Value = "Error"
$ws6.Range("G2").Value = 119802

$rng6 = $ws6.Range("A1:I2")
$rng6.Borders.LineStyle = 1
$rng6.Borders.Weight = 2
$rng6.BorderAround(1, -4138)

$hdr6 = $ws6.Range("A1:I1")
$hdr6.Font.Bold = $true
$hdr6.Font.Size = 14
$ws6.Rows(1).RowHeight = 18.75
$ws6.Rows(2).RowHeight = 15.75

$ws6.Range("G2").HorizontalAlignment = -4131
$ws6.Range("A1:I1").Select()

# ---------------------------------------------------------------------------
# Workbook-level: move active tab from AccountBalanceInformation(2nd) to
# StandingOrderAuth (6th, last sheet) and select it.
# ---------------------------------------------------------------------------
$ws6.Activate()

Write-Host "edit complete"
